$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BS32 section (rows 4-7): fill in integrator-call counts
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 37
$ws.Range("D6").Value = 86
$ws.Range("D7").Value = 195

# Update the averaging formula in G33 (BS54 section) with new run data
$ws.Range("G33").Formula = "=(191+2*193+203)/4"

# BS32 section row 52: fill in integrator-call count for last tolerance
$ws.Range("D52").Value = 193

# Restore the view to the top of the sheet and select D8
$ws.Range("D8").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
